# Se deja a medias lo de crandall, hace falta verificar suma de proyecciones
#
# This script updates the "Pol_rtas" polygon/traverse workbook:
#  - Renames survey points CT21 -> CD20, CT20 -> CD17, C20 -> D4, and
#    repoints one station (A8) from CT21 to D3 (mid-edit of the traverse).
#  - Updates a handful of computed projection / coordinate values on the
#    "Coordenadas" sheet.
#  - Updates the permissible angular error label and a few computed
#    parameters on the "Parametros Pol" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Proyecciones": rename survey station labels
# ---------------------------------------------------------------------
$wsProy = $wb.Worksheets.Item("Proyecciones")

$wsProy.Range("A2").Value = "CD20"
$wsProy.Range("B2").Value = "CD17"
$wsProy.Range("B4").Value = "CD20"
$wsProy.Range("A8").Value = "D3"
$wsProy.Range("B9").Value = "D4"

# ---------------------------------------------------------------------
# Sheet "Coordenadas": updated projection / coordinate correction values
# ---------------------------------------------------------------------
$wsCoord = $wb.Worksheets.Item("Coordenadas")

$wsCoord.Range("C3").Value = 0.005
$wsCoord.Range("D3").Value = -0.005

$wsCoord.Range("C5").Value = 0.007
$wsCoord.Range("D5").Value = -0.001
$wsCoord.Range("F5").Value = 2103.117

$wsCoord.Range("C7").Value = 0.002
$wsCoord.Range("D7").Value = -0.007
$wsCoord.Range("E7").Value = 1136.198
$wsCoord.Range("F7").Value = 2077.495

# ---------------------------------------------------------------------
# Sheet "Parametros Pol": updated permissible error + computed results
# ---------------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("Parametros Pol")

$wsParam.Range("B5").Value = "0° 0'30.0"

$wsParam.Range("D3").Value = -0.01399999999999935
$wsParam.Range("D4").Value = 0.01300000000000523
$wsParam.Range("D5").Value = 0.01910497317454588
$wsParam.Range("D6").Value = 11662.25139205389
